$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update workbook window tab ratio (bookViews / workbookView tabRatio 994 -> 988)
$wb.Windows.Item(1).TabRatio = 0.988

# Update header row
$ws.Range("B1").Value = "Tuotteen väri"
$ws.Range("C1").Value = "Malliston nimi"
$ws.Range("D1").Value = "poista"

# Row 2: hammer123 / musta / deluxe
$ws.Range("A2").Value = "hammer123"
$ws.Range("B2").Value = "musta"
$ws.Range("C2").Value = "deluxe"

# Row 3: helmet123 / sininen / sale
$ws.Range("A3").Value = "helmet123"
$ws.Range("B3").Value = "sininen"
$ws.Range("C3").Value = "sale"
$ws.Range("D3").Value = $null

# Row 4: helmet123 / (blank) / (blank) / X
$ws.Range("A4").Value = "helmet123"
$ws.Range("D4").Value = "X"

# Row 5: ski1 / valkoinen / winter
$ws.Range("A5").Value = "ski1"
$ws.Range("B5").Value = "valkoinen"
$ws.Range("C5").Value = "winter"

# Selection / active cell at B5
$ws.Range("B5").Select()
